# Apply updated dSF (column F) values as part of "repull data, push all data, mean calculation"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F3").Value  = 4
$ws.Range("F6").Value  = -5
$ws.Range("F8").Value  = -3
$ws.Range("F14").Value = -8
$ws.Range("F18").Value = -2
$ws.Range("F19").Value = -3
$ws.Range("F20").Value = -5
$ws.Range("F21").Value = -4
$ws.Range("F22").Value = 1
